$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 38.93002066666667
$ws.Range("H2").Value = 116.790062
$ws.Range("I2").Value = 0.7610372167397395
$ws.Range("J2").Value = 0.7610372167397393
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 6544.519052356432
$ws.Range("R2").Value = 58900.67147120788
$ws.Range("S2").Value = 0.2271073759897179
$ws.Range("T2").Value = 0.227107375989718
$ws.Range("G3").Value = 38.93002066666667
$ws.Range("H3").Value = 116.790062
$ws.Range("I3").Value = 0.7610372167397395
$ws.Range("J3").Value = 0.7610372167397393
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 6345.836123298871
$ws.Range("R3").Value = 57112.52510968984
$ws.Range("S3").Value = 0.2202126968985223
$ws.Range("T3").Value = 0.2202126968985222
$ws.Range("G4").Value = 38.93002066666667
$ws.Range("H4").Value = 116.790062
$ws.Range("I4").Value = 0.7610372167397395
$ws.Range("J4").Value = 0.7610372167397393
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 6462.131553432954
$ws.Range("R4").Value = 58159.18398089658
$ws.Range("S4").Value = 0.2242483716006743
$ws.Range("T4").Value = 0.2242483716006743
$ws.Range("G5").Value = 38.93002066666667
$ws.Range("H5").Value = 116.790062
$ws.Range("I5").Value = 0.7610372167397395
$ws.Range("J5").Value = 0.7610372167397393
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 2578.20813628251
$ws.Range("R5").Value = 23203.87322654259
$ws.Range("S5").Value = 0.0894687722508249
$ws.Range("T5").Value = 0.08946877225082489
$ws.Range("I6").Value = 0.1914142145281647
$ws.Range("J6").Value = 0.1914142145281647
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 1646.061383486611
$ws.Range("R6").Value = 14814.5524513795
$ws.Range("S6").Value = 0.05712149029301797
$ws.Range("T6").Value = 0.05712149029301797
$ws.Range("I7").Value = 0.1914142145281647
$ws.Range("J7").Value = 0.1914142145281647
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.05538735751522991
$ws.Range("T7").Value = 0.0553873575152299
$ws.Range("I8").Value = 0.1914142145281647
$ws.Range("J8").Value = 0.1914142145281647
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 1625.339481789167
$ws.Range("R8").Value = 14628.0553361025
$ws.Range("S8").Value = 0.05640240051997667
$ws.Range("T8").Value = 0.05640240051997666
$ws.Range("I9").Value = 0.1914142145281647
$ws.Range("J9").Value = 0.1914142145281647
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 648.4645881193612
$ws.Range("R9").Value = 5836.18129307425
$ws.Range("S9").Value = 0.02250296619994019
$ws.Range("T9").Value = 0.02250296619994019
$ws.Range("G10").Value = 1.794146
$ws.Range("H10").Value = 5.382438
$ws.Range("I10").Value = 0.03507349482179579
$ws.Range("J10").Value = 0.03507349482179579
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 301.6135742707907
$ws.Range("R10").Value = 2714.522168437116
$ws.Range("S10").Value = 0.01046657009743985
$ws.Range("T10").Value = 0.01046657009743985
$ws.Range("G11").Value = 1.794146
$ws.Range("H11").Value = 5.382438
$ws.Range("I11").Value = 0.03507349482179579
$ws.Range("J11").Value = 0.03507349482179579
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 292.4569856964074
$ws.Range("R11").Value = 2632.112871267666
$ws.Range("S11").Value = 0.01014881889410324
$ws.Range("T11").Value = 0.01014881889410323
$ws.Range("G12").Value = 1.794146
$ws.Range("H12").Value = 5.382438
$ws.Range("I12").Value = 0.03507349482179579
$ws.Range("J12").Value = 0.03507349482179579
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 297.81662787538
$ws.Range("R12").Value = 2680.34965087842
$ws.Range("S12").Value = 0.01033480876773223
$ws.Range("T12").Value = 0.01033480876773223
$ws.Range("G13").Value = 1.794146
$ws.Range("H13").Value = 5.382438
$ws.Range("I13").Value = 0.03507349482179579
$ws.Range("J13").Value = 0.03507349482179579
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 118.8204304972127
$ws.Range("R13").Value = 1069.383874474914
$ws.Range("S13").Value = 0.004123297062520487
$ws.Range("T13").Value = 0.004123297062520486
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.6381486666666666
$ws.Range("H14").Value = 1.914446
$ws.Range("I14").Value = 0.01247507391030006
$ws.Range("J14").Value = 0.01247507391030006
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 107.2790621663302
$ws.Range("R14").Value = 965.5115594969719
$ws.Range("S14").Value = 0.003722789423076182
$ws.Range("T14").Value = 0.003722789423076182
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.6381486666666666
$ws.Range("H15").Value = 1.914446
$ws.Range("I15").Value = 0.01247507391030006
$ws.Range("J15").Value = 0.01247507391030006
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 104.0222119490358
$ws.Range("R15").Value = 936.199907541322
$ws.Range("S15").Value = 0.003609770467684042
$ws.Range("T15").Value = 0.003609770467684042
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.6381486666666666
$ws.Range("H16").Value = 1.914446
$ws.Range("I16").Value = 0.01247507391030006
$ws.Range("J16").Value = 0.01247507391030006
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 105.9285498447933
$ws.Range("R16").Value = 953.3569486031399
$ws.Range("S16").Value = 0.003675924052659759
$ws.Range("T16").Value = 0.003675924052659759
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.6381486666666666
$ws.Range("H17").Value = 1.914446
$ws.Range("I17").Value = 0.01247507391030006
$ws.Range("J17").Value = 0.01247507391030006
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 42.26250221250422
$ws.Range("R17").Value = 380.362519912538
$ws.Range("S17").Value = 0.001466589966880082
$ws.Range("T17").Value = 0.001466589966880082
